# edit.ps1
# Applies 4 textual changes to the UDrive manual per the commit diff:
# 1) Merge "P" + "rocesador..." runs into a single run.
# 2) Split the "...cantidad de pixel que..." run, wrapping "pixel" in
#    gramStart/gramEnd proofErr markers (simulated grammar-check split).
# 3) Rewrite the SERIALIZACIÓN DE LOS DATOS paragraph: split "que" from the
#    following comma, italicize the three stream class names, and replace
#    the "etc...por medio de..." tail with the expanded ArrayList text.
# 4) Split "A continuación..." so "continuación" is wrapped in
#    gramStart/gramEnd proofErr markers.

$d = $word.ActiveDocument

function Get-ParagraphIndexContaining($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like $needle) {
            return $i
        }
    }
    return -1
}

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
$xmlFooter = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# Change 1: "P" + "rocesador: 2 GHz o superior" -> "Procesador: 2 GHz o superior"
# A same-text Find/Replace over the merged phrase collapses the two runs
# into a single run without disturbing anything else.
# ---------------------------------------------------------------------
$procRange = $d.Content
$procRange.Find.Execute("Procesador: 2 GHz o superior", $true, $false, $false, $false, $false, $true, 1, $false, "Procesador: 2 GHz o superior", 2)
Write-Output ("Change1 Found: " + $procRange.Find.Found)

# ---------------------------------------------------------------------
# Change 2: split "...cantidad de pixel que..." around "pixel"
# ---------------------------------------------------------------------
$idx2 = Get-ParagraphIndexContaining $d "*VENTAN INICIO DE VIAJES*"
Write-Output ("Change2 paragraph index: " + $idx2)
$p2 = $d.Paragraphs.Item($idx2)
$rng2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$xml2 = $xmlHeader + '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">VENTAN INICIO DE VIAJES: </w:t></w:r><w:r><w:t xml:space="preserve">Para la creación de esta ventana que es la más compleja de todas, se creó un ciclo el cual verificaría el valor del tipo de vehículo que se escogió (que está dentro de un vector previamente mencionado), para que de este modo se pudiera se coloque la imagen correspondiente al vehículo. Además, para que todo se colocara en orden , dependiendo de la posición, se sumó una cierta cantidad de </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>pixel</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> que distanciaría cada una de las 3 rutas de las otras. Dentro de todo esto, se colocan los datos de cada viaje, la gasolina, etc.</w:t></w:r>' + $xmlFooter
$rng2.InsertXML($xml2)
Write-Output "Change2 applied"

# ---------------------------------------------------------------------
# Change 3: rewrite the SERIALIZACIÓN DE LOS DATOS paragraph
# ---------------------------------------------------------------------
$idx3 = Get-ParagraphIndexContaining $d "*SERIALIZACIÓN DE LOS DATOS*"
Write-Output ("Change3 paragraph index: " + $idx3)
$p3 = $d.Paragraphs.Item($idx3)
$rng3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$xml3 = $xmlHeader + '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">SERIALIZACIÓN DE LOS DATOS: </w:t></w:r><w:r><w:t xml:space="preserve">Este programa requirió que los valores se serializaran para que no se pierdan al cerrar el programa. Por lo </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>que</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> para la serialización, se usaron </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006271E5"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>ObjectInputStream</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r><w:r w:rsidRPr="006271E5"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006271E5"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>FileInputStream</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r><w:r w:rsidRPr="006271E5"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006271E5"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>ObjectOutputStream</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> entre otros, p</w:t></w:r><w:r><w:t>ara guardar o leer el archivo que se generaría o guardaría</w:t></w:r><w:r><w:t>, se usó u</w:t></w:r><w:r><w:t xml:space="preserve">na matriz </w:t></w:r><w:r><w:t xml:space="preserve">de tipo objeto para guardar este tipo de información y un </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ArrayList</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, también de tipo objeto, para que se p</w:t></w:r><w:r><w:t>udiera guardar el historial.</w:t></w:r>' + $xmlFooter
$rng3.InsertXML($xml3)
Write-Output "Change3 applied"

# ---------------------------------------------------------------------
# Change 4: split "A continuación se dará..." around "continuación"
# ---------------------------------------------------------------------
$idx4 = Get-ParagraphIndexContaining $d "*A continuación se dará*"
Write-Output ("Change4 paragraph index: " + $idx4)
$p4 = $d.Paragraphs.Item($idx4)
$rng4 = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$xml4 = $xmlHeader + '<w:r><w:t xml:space="preserve">A </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>continuación</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> se dará una explicación general de lo que hace cada función:</w:t></w:r>' + $xmlFooter
$rng4.InsertXML($xml4)
Write-Output "Change4 applied"
